# capitalize first letter of every word in "Area of Expertise" (column B)
# and strip stray trailing newlines from a handful of "Contact Details" (column H) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Area of Expertise (column B): Title-case the category labels ---
$colB = $ws.Columns("B:B")
$colB.Replace("art-and-craft", "Art And Craft", 1)
$colB.Replace("creators", "Creators", 1)
$colB.Replace("dance", "Dance", 1)
$colB.Replace("photography", "Photography", 1)
$colB.Replace("poetry", "Poetry", 1)
$colB.Replace("theatre-drama", "Theatre Drama", 1)

# --- Contact Details (column H): drop the trailing newline characters ---
$ws.Range("H9").Value = "thekilljoycat@gmail.com"
$ws.Range("H17").Value = "postoncards@gmail.com"
$ws.Range("H37").Value = "artdesk.saara@gmail.comCommissions"
$ws.Range("H43").Value = "commissionsamirkhanpathanstudio@gmail.com"
$ws.Range("H49").Value = "ageisjustanumber62@gmail.com"
$ws.Range("H51").Value = "d@vasantiakhani2.0"
$ws.Range("H52").Value = "Collab-mailishadang@gmail.comTutorials"
$ws.Range("H55").Value = "reach.nishchay@gmail.com"
$ws.Range("H66").Value = "Sudhir@SudhirShivaram.com"
$ws.Range("H72").Value = "Kushagra@exifmedia.comCreative"
$ws.Range("H74").Value = "ajinkyakalbhor59@gmail.com"
$ws.Range("H77").Value = "Ulachhiramka@gmail.com"
$ws.Range("H78").Value = "snapwithpankaj@gmail.com"
$ws.Range("H91").Value = "subhashs.in@gmail.com"
$ws.Range("H108").Value = "kunalmalhotrakunal.malhotra2604@gmail.comLEARN"
$ws.Range("H109").Value = "roshanishah181@gmail.com"
$ws.Range("H127").Value = "kopal@tapeatale.com"
